$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("approvalDetails")
$ws.Cells.Item(1,1).Value = "test"
